$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("C18").Value = "['MEC-2NA-Fundição', -, -, -]"
$ws.Range("D18").Value = "[-, -, 'MEC-2NA-Fundição', -]"

# Row 19
$ws.Range("C19").Value = "['MEC-2NA-Fundição', -, -, -]"
$ws.Range("D19").Value = "[-, -, 'MEC-2NA-Fundição', -]"

# Row 20
$ws.Range("C20").Value = "-"

# Row 21
$ws.Range("C21").Value = "-"
